$words = @(
    "spüren",
    "proben",
    "wenden",
    "graben",
    "eignen",
    "duschen",
    "gelten",
    "messen",
    "achten",
    "lügen",
    "warnen",
    "mögen",
    "kosten",
    "zeichnen",
    "stammen",
    "fangen",
    "schmecken",
    "reisen",
    "heißen",
    "münzen",
    "kichern",
    "zählen",
    "altern",
    "heulen",
    "pfeifen",
    "jubeln",
    "gründen",
    "ruhen",
    "weinen",
    "feiern",
    "biegen",
    "zünden",
    "sprechen",
    "schlucken",
    "melden",
    "äußern",
    "werfen",
    "warten",
    "klettern",
    "seufzen",
    "platzen",
    "weichen",
    "schrecken",
    "brauchen",
    "hacken",
    "sterben",
    "führen",
    "rufen",
    "räumen",
    "wandern",
    "quälen",
    "fahren",
    "stehlen",
    "heilen",
    "tragen",
    "dienen",
    "parken",
    "fällen",
    "wundern",
    "kümmern",
    "backen",
    "zielen",
    "prügeln",
    "decken",
    "fischen",
    "rasen",
    "klingen",
    "leiden",
    "bergen",
    "filmen",
    "meinen",
    "spenden",
    "ärgern",
    "lohnen",
    "formen",
    "schultern",
    "bauen",
    "füllen",
    "boxen",
    "helfen",
    "grüßen",
    "kennen",
    "runden",
    "suchen",
    "schütteln",
    "passen",
    "knarren",
    "herrschen",
    "machen",
    "bilden",
    "betteln",
    "wohnen",
    "ändern",
    "sitzen",
    "ehren",
    "kehren",
    "liegen",
    "kratzen",
    "schreiten",
    "sichern",
    "malen",
    "leugnen",
    "zeigen",
    "husten",
    "trauen",
    "testen",
    "plaudern",
    "sprengen",
    "scheitern",
    "beißen",
    "albern",
    "werden",
    "streichen",
    "hören",
    "schützen",
    "planen",
    "heben",
    "sinken",
    "fließen",
    "greifen",
    "sperren",
    "fühlen",
    "lieben",
    "dauern",
    "fügen",
    "zögern",
    "wirken",
    "töten",
    "irren",
    "bitten",
    "wollen",
    "erben",
    "pressen",
    "fallen",
    "arten",
    "spannen",
    "flehen",
    "folgen",
    "drehen",
    "schulden",
    "starten",
    "saufen",
    "jagen",
    "streifen",
    "atmen",
    "freuen",
    "hauen",
    "pflanzen",
    "rauben",
    "treiben",
    "schenken",
    "lehnen",
    "stecken",
    "tollen",
    "küssen",
    "segeln",
    "enden",
    "lösen",
    "trennen",
    "dringen",
    "schwächen",
    "wachsen",
    "krachen",
    "geben",
    "ahnen",
    "mauern",
    "spielen",
    "wüten",
    "flüchten",
    "siegen",
    "spinnen",
    "locken",
    "schwingen",
    "sorgen",
    "bluten",
    "pissen",
    "bellen",
    "schwören",
    "drohen",
    "scheinen",
    "wehtun",
    "kaufen",
    "reizen",
    "lächeln",
    "fehlen",
    "liefern",
    "schwimmen",
    "lesen",
    "loben",
    "lockern",
    "tropfen",
    "schlafen"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
